$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 2682.3333
$ws.Range("I21").Value = 52
$ws.Range("K21").Value = 52
$ws.Range("M21").Value = 416

$ws.Range("H23").Value = 2682.3333
$ws.Range("I23").Value = 52
$ws.Range("K23").Value = 52
$ws.Range("M23").Value = 182

$ws.Range("H31").Value = 750
$ws.Range("J31").Value = 1000
$ws.Range("L31").Value = 3000
$ws.Range("N31").Value = -3460

$ws.Range("H40").Value = 1121.6428
$ws.Range("I40").Value = 1130.1
$ws.Range("J40").Value = 1100.5
$ws.Range("K40").Value = 1130.1
$ws.Range("L40").Value = 1100.5
$ws.Range("M40").Value = -955.0999999999999
$ws.Range("N40").Value = -1450.5

$ws.Range("H52").Value = 1920
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 1920
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 5760
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -6080

$ws.Range("H58").Value = 1115
$ws.Range("I58").Value = 172.5
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 517.5
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -367.5
$ws.Range("N58").Value = -9300

$ws.Range("H100").Value = 33336734
$ws.Range("I100").Value = 55558224
$ws.Range("K100").Value = 55558224
$ws.Range("M100").Value = -55557683

$ws.Range("H137").Value = 1785
$ws.Range("I137").Value = 1670
$ws.Range("J137").Value = 1900
$ws.Range("K137").Value = 5010
$ws.Range("L137").Value = 5700
$ws.Range("M137").Value = -2460
$ws.Range("N137").Value = -10800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 4634
$ws.Range("J46").Value = 3716.6667
$ws.Range("L46").Value = 3716.6667
$ws.Range("N46").Value = -4354.6667

$ws.Range("H61").Value = 1601.1111
$ws.Range("I61").Value = 1345.238
$ws.Range("J61").Value = 2496.6667
$ws.Range("K61").Value = 1345.238
$ws.Range("L61").Value = 2496.6667
$ws.Range("M61").Value = -1133.238
$ws.Range("N61").Value = -2920.6667

$ws.Range("H74").Value = 1385.6
$ws.Range("I74").Value = 1175.5
$ws.Range("J74").Value = 1625.7142
$ws.Range("K74").Value = 1175.5
$ws.Range("L74").Value = 1625.7142
$ws.Range("M74").Value = -301.5
$ws.Range("N74").Value = -3373.7142

$ws.Range("H77").Value = 1385.6
$ws.Range("I77").Value = 1175.5
$ws.Range("J77").Value = 1625.7142
$ws.Range("K77").Value = 5877.5
$ws.Range("L77").Value = 8128.571
$ws.Range("M77").Value = -1509.5
$ws.Range("N77").Value = -16864.571

$ws.Range("H136").Value = 1601.1111
$ws.Range("I136").Value = 1345.238
$ws.Range("J136").Value = 2496.6667
$ws.Range("K136").Value = 4035.714
$ws.Range("L136").Value = 7490.000100000001
$ws.Range("M136").Value = -1485.714
$ws.Range("N136").Value = -12590.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1743643.8
$ws.Range("I31").Value = 1931849.1
$ws.Range("J31").Value = 2744.5
$ws.Range("K31").Value = 1931849.1
$ws.Range("L31").Value = 2744.5
$ws.Range("M31").Value = -1931554.1
$ws.Range("N31").Value = -3334.5

$ws.Range("H34").Value = 1743643.8
$ws.Range("I34").Value = 1931849.1
$ws.Range("J34").Value = 2744.5
$ws.Range("K34").Value = 1931849.1
$ws.Range("L34").Value = 2744.5
$ws.Range("M34").Value = -1931647.1
$ws.Range("N34").Value = -3148.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1228
$ws.Range("I70").Value = 970.6667
$ws.Range("K70").Value = 2912.0001
$ws.Range("M70").Value = -2597.0001

$ws.Range("H73").Value = 1228
$ws.Range("I73").Value = 970.6667
$ws.Range("K73").Value = 2912.0001
$ws.Range("M73").Value = -1820.0001

$ws.Range("H76").Value = 4569.5625
$ws.Range("I76").Value = 1463
$ws.Range("J76").Value = 4776.6665
$ws.Range("K76").Value = 4389
$ws.Range("L76").Value = 14329.9995
$ws.Range("M76").Value = -4006
$ws.Range("N76").Value = -15095.9995

$ws.Range("H79").Value = 4569.5625
$ws.Range("I79").Value = 1463
$ws.Range("J79").Value = 4776.6665
$ws.Range("K79").Value = 4389
$ws.Range("L79").Value = 14329.9995
$ws.Range("M79").Value = -3063
$ws.Range("N79").Value = -16981.9995

$ws.Range("H80").Value = 2949.6667
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2949.6667
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 8849.000100000001
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -10721.0001

$ws.Range("H81").Value = 5279.2
$ws.Range("I81").Value = 145
$ws.Range("J81").Value = 5725.6523
$ws.Range("K81").Value = 435
$ws.Range("L81").Value = 17176.9569
$ws.Range("M81").Value = 688
$ws.Range("N81").Value = -19422.9569

$ws.Range("H82").Value = 1964.125
$ws.Range("I82").Value = 237.66667
$ws.Range("K82").Value = 713.00001
$ws.Range("M82").Value = -307.00001

$ws.Range("H83").Value = 2949.6667
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2949.6667
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 26547.0003
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -35907.0003

$ws.Range("H84").Value = 5279.2
$ws.Range("I84").Value = 145
$ws.Range("J84").Value = 5725.6523
$ws.Range("K84").Value = 1305
$ws.Range("L84").Value = 51530.8707
$ws.Range("M84").Value = 4311
$ws.Range("N84").Value = -62762.8707

$ws.Range("H85").Value = 1964.125
$ws.Range("I85").Value = 237.66667
$ws.Range("K85").Value = 713.00001
$ws.Range("M85").Value = 690.99999

$ws.Range("H92").Value = 533.3333
$ws.Range("I92").Value = 550
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 1650
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = -402
$ws.Range("N92").Value = -3996

$ws.Range("H93").Value = 5000
$ws.Range("J93").Value = 5000
$ws.Range("L93").Value = 15000
$ws.Range("N93").Value = -18744

$ws.Range("H107").Value = 333.18182
$ws.Range("I107").Value = 480
$ws.Range("J107").Value = 318.5
$ws.Range("K107").Value = 1440
$ws.Range("L107").Value = 955.5
$ws.Range("M107").Value = 480
$ws.Range("N107").Value = -4795.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 7000
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws.Range("H136").Value = 26715.572
$ws.Range("I136").Value = 21420.8
$ws.Range("J136").Value = 39952.5
$ws.Range("K136").Value = 64262.39999999999
$ws.Range("L136").Value = 119857.5
$ws.Range("M136").Value = -61712.39999999999
$ws.Range("N136").Value = -124957.5
